# Weekly price update for "Hortaliza, Macroferia Regional de Talca - Poroto verde".
#
# A new price record (date 2022-12-02 / serial 44897) is inserted as row 95,
# pushing the existing rows 95:185 down to 96:186 (dimension grows from
# A1:R185 to A1:R186). The new row copies the qualitative fields (variety,
# quality, unit of sale, origin, kg-per-unit, classification) from the row
# it displaces, but carries its own volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 95; Excel shifts rows 95:185 -> 96:186 automatically.
$ws.Rows("95:95").Insert()

# Populate the newly inserted row 95 with the new record.
$ws.Range("A95").Value = 5
$ws.Range("B95").Value = "Macroferia Regional de Talca"
$ws.Range("C95").Value = "Maule"
$ws.Range("D95").Value = 44897
$ws.Range("E95").Value = 7
$ws.Range("F95").Value = 100112031
$ws.Range("G95").Value = "Poroto verde"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 50
$ws.Range("K95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("M95").Value = 35000
$ws.Range("N95").Value = "$/saco 25 kilos"
$ws.Range("O95").Value = "Región del Maule"
$ws.Range("P95").Value = 1400
$ws.Range("Q95").Value = 25
$ws.Range("R95").Value = "Hortaliza"
